# daily auto push: 2026-02-07 06:59 UTC
# Insert a new data row for 2026/02/07 (hour 14) just before the
# 2026/12/29 block, shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 795; existing rows 795.. shift down to 796..
$ws.Rows.Item(795).Insert()

# Populate the newly inserted row with the new entry's values.
# The date column holds plain text (not a real date), so force text
# formatting before assignment to stop Excel from auto-converting the
# "yyyy/mm/dd"-looking string into a date serial, then restore the
# default "Normal" style so the cell's style matches its siblings.
$ws.Cells.Item(795, 1).NumberFormat = "@"
$ws.Cells.Item(795, 1).Value = "2026/02/07"
$ws.Cells.Item(795, 1).Style = "Normal"

$ws.Cells.Item(795, 2).Value = "土"
$ws.Cells.Item(795, 3).Value = 14
$ws.Cells.Item(795, 4).Value = 201
